$wb = $excel.ActiveWorkbook

# Row update 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2804.9443
$ws.Range("I116").Value = 2967.6667
$ws.Range("J116").Value = 2642.2222
$ws.Range("K116").Value = 2967.6667
$ws.Range("L116").Value = 2642.2222
$ws.Range("M116").Value = 474.3332999999998
$ws.Range("N116").Value = -9526.2222

# Row update 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2535.6562
$ws.Range("I132").Value = 1994.0385
$ws.Range("J132").Value = 4882.6665
$ws.Range("K132").Value = 5982.1155
$ws.Range("L132").Value = 14647.9995
$ws.Range("M132").Value = -3452.1155
$ws.Range("N132").Value = -19707.9995

# Row update 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3357.566
$ws.Range("I137").Value = 1038.0333
$ws.Range("J137").Value = 6383.0435
$ws.Range("K137").Value = 3114.0999
$ws.Range("L137").Value = 19149.1305
$ws.Range("M137").Value = -564.0999000000002
$ws.Range("N137").Value = -24249.1305

# Row update 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1560.7167
$ws.Range("I138").Value = 906.44684
$ws.Range("J138").Value = 3926.1538
$ws.Range("K138").Value = 2719.34052
$ws.Range("L138").Value = 11778.4614
$ws.Range("M138").Value = 2420.65948
$ws.Range("N138").Value = -22058.4614

# Row update 4: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 125863.75
$ws.Range("I2").Value = 333970.34
$ws.Range("J2").Value = 999.8
$ws.Range("K2").Value = 333970.34
$ws.Range("L2").Value = 999.8
$ws.Range("M2").Value = -333857.34
$ws.Range("N2").Value = -1225.8

# Row update 5: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2292.5
$ws.Range("I63").Value = 2337.2727
$ws.Range("J63").Value = 1800
$ws.Range("K63").Value = 2337.2727
$ws.Range("L63").Value = 1800
$ws.Range("M63").Value = -1651.2727
$ws.Range("N63").Value = -3172

# Row update 6: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2292.5
$ws.Range("I66").Value = 2337.2727
$ws.Range("J66").Value = 1800
$ws.Range("K66").Value = 11686.3635
$ws.Range("L66").Value = 9000
$ws.Range("M66").Value = -8254.363499999999
$ws.Range("N66").Value = -15864

# Row update 7: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1101.4286
$ws.Range("I102").Value = 951.6667
$ws.Range("K102").Value = 951.6667
$ws.Range("M102").Value = 670.3333

# Row update 8: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 125863.75
$ws.Range("I116").Value = 333970.34
$ws.Range("J116").Value = 999.8
$ws.Range("K116").Value = 333970.34
$ws.Range("L116").Value = 999.8
$ws.Range("M116").Value = -331676.34
$ws.Range("N116").Value = -5587.8

# Row update 9: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 42079.25
$ws.Range("J134").Value = 42079.25
$ws.Range("L134").Value = 42079.25
$ws.Range("N134").Value = -52219.25

# Row update 10: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 125863.75
$ws.Range("I3").Value = 333970.34
$ws.Range("J3").Value = 999.8
$ws.Range("K3").Value = 333970.34
$ws.Range("L3").Value = 999.8
$ws.Range("M3").Value = -333856.34
$ws.Range("N3").Value = -1227.8

# Row update 11: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row update 12: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 886.15
$ws.Range("I94").Value = 721.0833
$ws.Range("J94").Value = 1133.75
$ws.Range("K94").Value = 721.0833
$ws.Range("L94").Value = 1133.75
$ws.Range("M94").Value = -270.0833
$ws.Range("N94").Value = -2035.75

# Row update 13: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# Row update 14: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34394.69
$ws.Range("I31").Value = 3162.3953
$ws.Range("J31").Value = 88114.24000000001
$ws.Range("K31").Value = 3162.3953
$ws.Range("L31").Value = 88114.24000000001
$ws.Range("M31").Value = -2867.3953
$ws.Range("N31").Value = -88704.24000000001

# Row update 15: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 34394.69
$ws.Range("I34").Value = 3162.3953
$ws.Range("J34").Value = 88114.24000000001
$ws.Range("K34").Value = 3162.3953
$ws.Range("L34").Value = 88114.24000000001
$ws.Range("M34").Value = -2960.3953
$ws.Range("N34").Value = -88518.24000000001

# Row update 16: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1386.8846
$ws.Range("I58").Value = 1418.1578
$ws.Range("J58").Value = 1302
$ws.Range("K58").Value = 1418.1578
$ws.Range("L58").Value = 1302
$ws.Range("M58").Value = -1215.1578
$ws.Range("N58").Value = -1708

# Row update 17: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2553.4546
$ws.Range("I99").Value = 1759.25
$ws.Range("J99").Value = 4671.3335
$ws.Range("K99").Value = 1759.25
$ws.Range("L99").Value = 4671.3335
$ws.Range("M99").Value = -261.25
$ws.Range("N99").Value = -7667.3335

# Row update 18: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2553.4546
$ws.Range("I126").Value = 1759.25
$ws.Range("J126").Value = 4671.3335
$ws.Range("K126").Value = 5277.75
$ws.Range("L126").Value = 14014.0005
$ws.Range("M126").Value = -2807.75
$ws.Range("N126").Value = -18954.0005

# Row update 19: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2873.0168
$ws.Range("I134").Value = 3319.5652
$ws.Range("J134").Value = 1292.9231
$ws.Range("K134").Value = 9958.695599999999
$ws.Range("L134").Value = 3878.7693
$ws.Range("M134").Value = -7423.695599999999
$ws.Range("N134").Value = -8948.7693

# Row update 20: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1386.8846
$ws.Range("I136").Value = 1418.1578
$ws.Range("J136").Value = 1302
$ws.Range("K136").Value = 4254.4734
$ws.Range("L136").Value = 3906
$ws.Range("M136").Value = -1704.4734
$ws.Range("N136").Value = -9006

# Row update 21: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 37058
$ws.Range("J141").Value = 37058
$ws.Range("L141").Value = 37058
$ws.Range("N141").Value = -47418

# Row update 22: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 33534964
$ws.Range("J22").Value = 37038850
$ws.Range("L22").Value = 111116550
$ws.Range("N22").Value = -111116888

# Row update 23: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 33534964
$ws.Range("J27").Value = 37038850
$ws.Range("L27").Value = 111116550
$ws.Range("N27").Value = -111116754

# Row update 24: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 83334340
$ws.Range("J58").Value = 166667680
$ws.Range("L58").Value = 500003040
$ws.Range("N58").Value = -500003296

# Row update 25: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 66667250
$ws.Range("I64").Value = 726
$ws.Range("J64").Value = 333333340
$ws.Range("K64").Value = 2178
$ws.Range("L64").Value = 1000000020
$ws.Range("M64").Value = -1908
$ws.Range("N64").Value = -1000000560

# Row update 26: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 66667250
$ws.Range("I67").Value = 726
$ws.Range("J67").Value = 333333340
$ws.Range("K67").Value = 2178
$ws.Range("L67").Value = 1000000020
$ws.Range("M67").Value = -1242
$ws.Range("N67").Value = -1000001892

# Row update 27: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 292.13333
$ws.Range("I107").Value = 267
$ws.Range("J107").Value = 301.27274
$ws.Range("K107").Value = 801
$ws.Range("L107").Value = 903.81822
$ws.Range("M107").Value = 1119
$ws.Range("N107").Value = -4743.81822

# Row update 28: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 903.0909
$ws.Range("I132").Value = 787.1111
$ws.Range("J132").Value = 1425
$ws.Range("K132").Value = 7083.9999
$ws.Range("L132").Value = 12825
$ws.Range("M132").Value = -4553.9999
$ws.Range("N132").Value = -17885

# Row update 29: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 22.166666
$ws.Range("J2").Value = 36.75
$ws.Range("K2").Value = 22.166666
$ws.Range("L2").Value = 36.75
$ws.Range("M2").Value = 90.83333400000001
$ws.Range("N2").Value = -262.75

# Row update 30: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5828.5454
$ws.Range("I122").Value = 6000.2856
$ws.Range("J122").Value = 2222
$ws.Range("K122").Value = 18000.8568
$ws.Range("L122").Value = 6666
$ws.Range("M122").Value = -15550.8568
$ws.Range("N122").Value = -11566

# Row update 31: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4529.5835
$ws.Range("I132").Value = 7073.3184
$ws.Range("J132").Value = 2377.1924
$ws.Range("K132").Value = 21219.9552
$ws.Range("L132").Value = 7131.5772
$ws.Range("M132").Value = -18689.9552
$ws.Range("N132").Value = -12191.5772

# Row update 32: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21430796
$ws.Range("I122").Value = 22223878
$ws.Range("J122").Value = 20411118
$ws.Range("K122").Value = 66671634
$ws.Range("L122").Value = 61233354
$ws.Range("M122").Value = -66669184
$ws.Range("N122").Value = -61238254
